$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) from _old/_new to _FV2404/_FV2410 ---
$oldHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i] + "_FV2404"
}

$newHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i] + "_FV2410"
}

# --- Add table over the whole data range, using header row for column names ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
